# Re-bin the age_group column (B) from 18-year-wide, 6-group bins into
# 10-year-wide bins with a new open-ended "78-" top bin, matching the
# new one-hot-encode-ready pivot grouping used in the notebook.
#
# Old bins: 18-35, 36-45, 46-55, 56-65, 66-75, 76-
# New bins: 18-27, 28-37, 38-47, 48-57, 58-67, 68-77, 78-

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row stays the same ("age" / "age_group"), only the group text
# for each data row (2..99, ages 18..115) is recomputed.
for ($row = 2; $row -le 99; $row++) {
    $age = $ws.Cells.Item($row, 1).Value()

    if ($age -le 27) {
        $group = "18-27"
    } elseif ($age -le 37) {
        $group = "28-37"
    } elseif ($age -le 47) {
        $group = "38-47"
    } elseif ($age -le 57) {
        $group = "48-57"
    } elseif ($age -le 67) {
        $group = "58-67"
    } elseif ($age -le 77) {
        $group = "68-77"
    } else {
        $group = "78-"
    }

    $ws.Cells.Item($row, 2).Value = $group
}

# Reflect the reviewed selection (the new open-ended "78-" bucket, rows 62:99)
$ws.Range("B62:B99").Select()
